$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header row: new "Status" column F1 (copy the header formatting from
# an existing header cell so font/border/alignment all match)
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 6).Value = "Status"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Row 2 (existing data row): fix Scalpel Type set ordering + add Status
# ---------------------------------------------------------------------
$ws.Cells.Item(2, 5).Value = "{'any', 'list', 'List[any]'}"
$ws.Cells.Item(2, 6).Value = "Neutral"

# ---------------------------------------------------------------------
# Row 3 (existing data row): PyType Type becomes "List", add Status
# ---------------------------------------------------------------------
$ws.Cells.Item(3, 4).Value = "List"
$ws.Cells.Item(3, 5).Value = "any"
$ws.Cells.Item(3, 6).Value = "Loss"

# ---------------------------------------------------------------------
# Row 4: used to be the "PyType Total / Scalpel Total" summary row,
# now becomes a third comparison data row
# ---------------------------------------------------------------------
$ws.Cells.Item(4, 1).Value = "minimaxir__big-list-of-naughty-strings"
$ws.Cells.Item(4, 2).Value = "__init__.py"
$ws.Cells.Item(4, 3).Value = "naughty_strings"
$ws.Cells.Item(4, 4).Value = "Any"
$ws.Cells.Item(4, 5).Value = "any"
$ws.Cells.Item(4, 6).Value = "Neutral"

# ---------------------------------------------------------------------
# Row 5: used to be the "Accuracy" row, now becomes the comparisons /
# wins summary row
# ---------------------------------------------------------------------
$ws.Cells.Item(5, 1).Value = "Total comparisons:"
$ws.Cells.Item(5, 2).Value = 3
$ws.Cells.Item(5, 3).Value = "PyType Wins:"
$ws.Cells.Item(5, 4).Value = 1
$ws.Cells.Item(5, 5).Value = "Scalpel Wins:"
$ws.Cells.Item(5, 6).Value = 0

# ---------------------------------------------------------------------
# Row 6 (new): accuracy-over-pytype row
# ---------------------------------------------------------------------
$ws.Cells.Item(6, 5).Value = "Accuracy over PyType"
$ws.Cells.Item(6, 6).Value = 0

# ---------------------------------------------------------------------
# Fill colors: white background on the whole A2:F6 data block, then
# override the Status cells with orange (Neutral) / red (Loss)
# ---------------------------------------------------------------------
$ws.Range("A2:F6").Interior.Color = 0x00FFFFFF

$ws.Range("F2").Interior.Color = 0x0000A5FF
$ws.Range("F3").Interior.Color = 0x000000FF
$ws.Range("F4").Interior.Color = 0x0000A5FF

Write-Host "edit applied"
